$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "beihai_adj_gaoqiazhen" row entirely (row 9), shifting all
# subsequent rows up by one. This matches the commit message:
# "removed beihi adj it was too much of a problem"
$ws.Rows.Item(9).Delete()
